$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New label next to the "conv_sed_mask(:)*dum_sed(:)" header row
$ws.Range("C7").Value = "THESE ARE GENIE VALUES (in cm3/(cm2*yr))"

# Make sure the new rows below the existing table use the sheet's normal row height
25..30 | ForEach-Object { $ws.Rows.Item($_).RowHeight = 12.1 }

# New block of conversion factors (mol -> cm3) appended below the existing table
$ws.Range("E26").Value = "Conversion factors mol " + [char]0x2192 + " cm3 (from GENIE)"
$ws.Range("E26").Font.Bold = $true

$ws.Range("E27").Value = "conv_POC_mol_cm3"
$ws.Range("F27").Formula = "=12"

$ws.Range("E28").Value = "conv_cal_mol_cm3"
$ws.Range("F28").Formula = "=100/2.7"

$ws.Range("E29").Value = "conv_det_mol_cm3"
$ws.Range("F29").Formula = "=60/3"

$ws.Range("E30").Value = "conv_ash_mol_cm3"

$ws.Range("I8").Select() | Out-Null
